$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Investment Type *" column (F) is removed entirely; everything to its
# right (Notes, Type, Folio No) shifts one column to the left.
$ws.Columns.Item(6).Delete()

# New trailing columns for category / sub-category classification of each
# valuation entry.
$ws.Range("I1").Value = "Sector"
$ws.Range("J1").Value = "Category"
$ws.Range("K1").Value = "Sub Category"
$ws.Range("L1").Value = "Startup"
$ws.Range("M1").Value = "Investment Domicile"

$ws.Range("J2:J9").Value = "Unlisted"
$ws.Range("K2:K9").Value = "Equity"
$ws.Range("L2:L9").Value = "Yes"
$ws.Range("M2:M9").Value = "Domestic"

$ws.Range("J3:J9").Select()
